# Update the "non-primitives" spec table to reflect the new table/range/formula
# examples and re-indented struct/class example, and adjust row heights / column
# widths / selection to match the refreshed layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update example cells in column C ---

$ws.Range("C3").Value = "range data = (1:5:1, 1:5:1);`nrange rng = (1:5:1);`nrange col = (2, 1:5:1);"

$ws.Range("C5").Value = "table tab1;"

$ws.Range("C6").Value = "formula f1 = SUM(1,1:10:2) / 10;"

$ws.Range("C7").Value = "class example {`n    int id;`n    string name;`n    double price;`n};`nclass example = {id: 101, name: `"widget`", price: 19.99}"

# --- Row heights ---

$ws.Rows.Item(1).RowHeight = 21
$ws.Rows.Item(2).RowHeight = 25.5
$ws.Rows.Item(3).RowHeight = 51
$ws.Rows.Item(5).RowHeight = 37
$ws.Rows.Item(6).RowHeight = 27.5
$ws.Rows.Item(7).RowHeight = 103.5

# --- Column widths (columns B and C made wider) ---

$ws.Columns.Item(2).ColumnWidth = 35.335
$ws.Columns.Item(3).ColumnWidth = 29.83

# --- Selection moves to E4 ---

$ws.Range("E4").Select()
